$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5: new runs (B) and matches played (D) values ---
$ws.Range("B2").Value = 693
$ws.Range("D2").Value = 10

$ws.Range("B3").Value = 385
$ws.Range("D3").Value = 10

$ws.Range("B4").Value = 578
$ws.Range("D4").Value = 10

$ws.Range("B5").Value = 473
$ws.Range("D5").Value = 10

# --- Add new rows 6-10 ---
# Column A (names) filled first, in order
$ws.Range("A6").Value = "WATSON"
$ws.Range("A7").Value = "SMITH"
$ws.Range("A8").Value = "GAYLE"
$ws.Range("A9").Value = "SACHIN"
$ws.Range("A10").Value = "DRAVID"

# Column B (runs)
$ws.Range("B6").Value = 398
$ws.Range("B7").Value = 659
$ws.Range("B8").Value = 603
$ws.Range("B9").Value = 690
$ws.Range("B10").Value = 573

# Column C (teams) filled after names
$ws.Range("C6").Value = "RR"
$ws.Range("C7").Value = "RR"
$ws.Range("C8").Value = "RCB"
$ws.Range("C9").Value = "MI"
$ws.Range("C10").Value = "DC"

# Column D (matches played)
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 10
$ws.Range("D9").Value = 10
$ws.Range("D10").Value = 10

# --- Average formula (column E) for all data rows 2-10 ---
$ws.Range("E2").Formula = "=ROUNDDOWN(B2/D2,2)"
$ws.Range("E3").Formula = "=ROUNDDOWN(B3/D3,2)"
$ws.Range("E4").Formula = "=ROUNDDOWN(B4/D4,2)"
$ws.Range("E5").Formula = "=ROUNDDOWN(B5/D5,2)"
$ws.Range("E6").Formula = "=ROUNDDOWN(B6/D6,2)"
$ws.Range("E7").Formula = "=ROUNDDOWN(B7/D7,2)"
$ws.Range("E8").Formula = "=ROUNDDOWN(B8/D8,2)"
$ws.Range("E9").Formula = "=ROUNDDOWN(B9/D9,2)"
$ws.Range("E10").Formula = "=ROUNDDOWN(B10/D10,2)"

# --- Update the selected range to match the new data extent ---
$ws.Range("E2:E10").Select() | Out-Null
